$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Date property value (row 8, column B)
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# Insert a new row for "Jurisdiction" (blank value) right after "Contact" (row 10)
# and before "Description" (old row 11), pushing everything below down by one row.
$ws.Rows.Item(11).Insert()

# Match the formatting/style used by the other data rows.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
